$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) Table cell "250000 images/day": drop the stray _GoBack bookmark that sits
#    between the "images" and "/day" runs. A self-replace via Find (wrap-all)
#    over a range spanning the bookmark strips the (contentless) bookmark.
# ---------------------------------------------------------------------------
$rngImages = $d.Content
$rngImages.Find.Execute("images/day", $true, $false, $false, $false, $false, $true, 1, $false, "images/day", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Measurements on a dual core MacOs machine[3] showed ..." paragraph:
#    split into several runs, add a proofErr gramStart/gramEnd pair around
#    "machine[" , and split the superscript "[3]" into "[" / "3]".
# ---------------------------------------------------------------------------
$rngMeasurements = $d.Content
$foundMeasurements = $rngMeasurements.Find.Execute("Measurements on a dual core MacOs machine", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMeasurements) {
    $measurementsPara = $rngMeasurements.Paragraphs(1).Range

    $measurementsBody = @'
<w:p><w:r><w:t xml:space="preserve">Measurements on </w:t></w:r><w:r><w:t xml:space="preserve">small sample of files, on a </w:t></w:r><w:r><w:t xml:space="preserve">dual core MacOs </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>machine</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>3]</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">showed </w:t></w:r><w:r><w:t>batch building process files rates of around 10 files/sec. Actual throughput increases for large batches (more files /sec).</w:t></w:r></w:p>
'@

    $measurementsXml = $pkgHeader + "<w:body>" + $measurementsBody + "</w:body>" + $pkgFooter
    $measurementsPara.InsertXML($measurementsXml) | Out-Null
}

# ---------------------------------------------------------------------------
# 3) "At this rate, the batch building steps would consume 6.4 days to
#    process." paragraph: change "6.4" -> "15" (as separate runs), then add a
#    brand-new paragraph describing the second test, ending with the _GoBack
#    bookmark (moved here from the table above).
# ---------------------------------------------------------------------------
$rngRate = $d.Content
$foundRate = $rngRate.Find.Execute("At this rate, the batch building steps would consume 6.4 days to process", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundRate) {
    $ratePara = $rngRate.Paragraphs(1).Range

    $rateAndSecondTestBody = @'
<w:p><w:r><w:t xml:space="preserve">At this rate, the batch building steps would consume </w:t></w:r><w:r><w:t>15</w:t></w:r><w:r><w:t xml:space="preserve"> days to process</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">A second </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>test  was</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> done processing 13 batches</w:t></w:r><w:r><w:t xml:space="preserve"> simultaneously. See sheet &#8216;batches113&#8217; in [1</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">] </w:t></w:r><w:r><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> This tripled the throughput, giving an estimated batch building </w:t></w:r><w:r><w:t xml:space="preserve">elapsed </w:t></w:r><w:r><w:t xml:space="preserve">time of </w:t></w:r><w:r><w:t>5  days.</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>
'@

    $rateAndSecondTestXml = $pkgHeader + "<w:body>" + $rateAndSecondTestBody + "</w:body>" + $pkgFooter
    $ratePara.InsertXML($rateAndSecondTestXml) | Out-Null
}

Write-Output "done"
